$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "09/04/2025"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = 0.1190458572798798
$ws.Range("C3").Value = 0.8809541427201202
